$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Headers: BTec_Logo-Orange picture, name="image2.jpg" -> name="image1.jpg" ---
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hf = $sec.Headers($i)
    if ($hf.Exists) {
        $xml = $hf.Range.WordOpenXML
        if ($xml -like "*BTec_Logo-Orange*") {
            $newXml = $xml.Replace('name="image2.jpg"', 'name="image1.jpg"')
            if ($newXml -ne $xml) {
                $hf.Range.WordOpenXML = $newXml
            }
        }
    }
}

# --- Footers: PearsonLogo picture, name="image1.png" -> name="image2.png" ---
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $hf = $sec.Footers($i)
    if ($hf.Exists) {
        $xml = $hf.Range.WordOpenXML
        if ($xml -like "*PearsonLogo*") {
            $newXml = $xml.Replace('name="image1.png"', 'name="image2.png"')
            if ($newXml -ne $xml) {
                $hf.Range.WordOpenXML = $newXml
            }
        }
    }
}

Write-Host "done"
